$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the per-fold metric values (re-run of the evaluation) in C2:L9
$ws.Range("C2").Value = 0.922289823008849
$ws.Range("D2").Value = 0.00785222933198126
$ws.Range("E2").Value = 0.97469315859758
$ws.Range("F2").Value = 0.00397066710546851
$ws.Range("G2").Value = 0.899171988388713
$ws.Range("H2").Value = 0.0150722697451505
$ws.Range("I2").Value = 0.926157215264941
$ws.Range("J2").Value = 0.0147247588875365
$ws.Range("K2").Value = 0.912309190821805
$ws.Range("L2").Value = 0.00907018295931141

$ws.Range("C3").Value = 0.9125
$ws.Range("D3").Value = 0.00865310146309931
$ws.Range("E3").Value = 0.96881585228145
$ws.Range("F3").Value = 0.00509310029991419
$ws.Range("G3").Value = 0.886934163270113
$ws.Range("H3").Value = 0.0148166042464536
$ws.Range("I3").Value = 0.914759603352811
$ws.Range("J3").Value = 0.0123283408584579
$ws.Range("K3").Value = 0.90053396226915
$ws.Range("L3").Value = 0.01001378553969

$ws.Range("C4").Value = 0.889712389380531
$ws.Range("D4").Value = 0.00920936145864612
$ws.Range("E4").Value = 0.888277185215031
$ws.Range("F4").Value = 0.00953623468885448
$ws.Range("G4").Value = 0.878442216496897
$ws.Range("H4").Value = 0.0218637165376288
$ws.Range("I4").Value = 0.866195068918889
$ws.Range("J4").Value = 0.0172418587933654
$ws.Range("K4").Value = 0.87198335763133
$ws.Range("L4").Value = 0.0114864638213478

$ws.Range("C5").Value = 0.882134955752212
$ws.Range("D5").Value = 0.0117926528171889
$ws.Range("E5").Value = 0.944415825243178
$ws.Range("F5").Value = 0.0079530702922135
$ws.Range("G5").Value = 0.851928945783977
$ws.Range("H5").Value = 0.0184024679792692
$ws.Range("I5").Value = 0.884161832701991
$ws.Range("J5").Value = 0.0171116276412259
$ws.Range("K5").Value = 0.867585538444625
$ws.Range("L5").Value = 0.0134749866429776

$ws.Range("C6").Value = 0.847898230088495
$ws.Range("D6").Value = 0.00752293722828288
$ws.Range("E6").Value = 0.920003446209121
$ws.Range("F6").Value = 0.00498620617248437
$ws.Range("G6").Value = 0.859953069513008
$ws.Range("H6").Value = 0.0119091954267867
$ws.Range("I6").Value = 0.78170764724566
$ws.Range("J6").Value = 0.0143494541724472
$ws.Range("K6").Value = 0.818868390845953
$ws.Range("L6").Value = 0.010089237597173

$ws.Range("C7").Value = 0.900940265486725
$ws.Range("D7").Value = 0.00960209488255584
$ws.Range("E7").Value = 0.96146874641498
$ws.Range("F7").Value = 0.00599802391784958
$ws.Range("G7").Value = 0.881924238810167
$ws.Range("H7").Value = 0.0183752295559158
$ws.Range("I7").Value = 0.892680972963396
$ws.Range("J7").Value = 0.0192351527533963
$ws.Range("K7").Value = 0.887045371822335
$ws.Range("L7").Value = 0.0124628301196515

$ws.Range("C8").Value = 0.882411504424778
$ws.Range("D8").Value = 0.00928215340522043
$ws.Range("E8").Value = 0.888254689371312
$ws.Range("F8").Value = 0.00842262854037781
$ws.Range("G8").Value = 0.869494737740491
$ws.Range("H8").Value = 0.0177187825987622
$ws.Range("I8").Value = 0.859510288983465
$ws.Range("J8").Value = 0.017339382279041
$ws.Range("K8").Value = 0.86427720641071
$ws.Range("L8").Value = 0.0116574193615812

$ws.Range("C9").Value = 0.909181415929203
$ws.Range("D9").Value = 0.00793766955792487
$ws.Range("E9").Value = 0.965534702361888
$ws.Range("F9").Value = 0.00555788500649627
$ws.Range("G9").Value = 0.88608076599905
$ws.Range("H9").Value = 0.0167716903788911
$ws.Range("I9").Value = 0.911783788177325
$ws.Range("J9").Value = 0.0149683066006302
$ws.Range("K9").Value = 0.898562562252537
$ws.Range("L9").Value = 0.00934213904522203

# The refreshed raw metrics no longer carry the 0.000 number format
$ws.Range("C2:L9").Style = "Normal"

# Selection moved by the author while reviewing the updated figures
$ws.Range("N5").Select() | Out-Null
